$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 corresponds to de0447c2-... file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: row 3 corresponds to de0447c2-... file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-01-27 08:16:36"

# --- de-de sheet: row 3 corresponds to de0447c2-... file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-01-27 08:16:48"
